$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column C ("Runmode"), shifting County/FirstName/Age/... right.
$ws.Columns("C:C").Insert()

# 2. Remove the now-obsolete bottom rows (10-16) that only ever held blank placeholders.
$ws.Rows("10:16").Delete()

# 3. Populate the header for the new Runmode column.
$ws.Range("C1").Value = "Runmode"
$ws.Range("C1").Borders.LineStyle = -4142

# 4. Populate Runmode values for the four existing data rows.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "Y"

# 5. Fill in the rest of row 3/4/5 with the same reference data already present in row 2
#    (County/FirstName/Age/Gender/CountOfPeopleLive/CountOfUnder19), matching row 2's pattern.
$ws.Range("D2:I2").Copy()
$ws.Range("D3:I3").Value = $ws.Range("D2:I2").Value()
$ws.Range("D4:I4").Value = $ws.Range("D2:I2").Value()
$ws.Range("D5:I5").Value = $ws.Range("D2:I2").Value()

# 6. Leave the selection where the author left it.
$ws.Range("D15").Select()
